$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: update title (D) and link (E)
$ws.Range("D28").Value = "MIT 6.800/6.843 Robotics Manipulation :: Introduction"
$ws.Range("E28").Value = "https://ropiens.tistory.com/191"

# Row 32: update title (D) and link (E)
$ws.Range("D32").Value = "All about Feature Scaling"
$ws.Range("E32").Value = "https://dodonam.tistory.com/376"

# Row 51: update title (D) and link (E)
$ws.Range("D51").Value = "[Mac] 맥북에서 한영 전환하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/1306"
